# payment + send invoice
# Adds a new Doctor Availability row (row 2) for D001, with the newly
# confirmed/sent invoice slot date-time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "D001"
$ws.Range("B2").Value = "2029-03-08 10:30"
